# Quarterly sheet ("Quarterly" / sheet4.xml): split the combined
# "<Mon> <YY> Q<n>" label column into three separate columns
# (Year / Month / Quarter) by inserting two new columns after column A
# and re-populating the period columns for each data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Quarterly")

# Insert two new blank columns at B:C - everything from the old column B
# onward (the numeric figures) shifts right by two columns (B->D, C->E,
# ... U->W) while column A (the title column) stays put.
$ws.Range("B1:C1").EntireColumn.Insert()

# New header labels for the inserted columns.
$ws.Range("B1").Value = "Month"
$ws.Range("C1").Value = "Quarter"

# The Year/Month columns hold text that looks numeric ("2023", "09", ...)
# so format them as Text first - otherwise Excel would coerce them to
# numbers (dropping the leading zero on the month code).
$ws.Range("A2:B6").NumberFormat = "@"

# Re-populate column A (year), and the newly-inserted B (month) / C
# (quarter) columns for each of the five quarterly data rows, replacing
# the old single "Mon YY Qn" label.
$ws.Range("A2").Value = "2023"
$ws.Range("B2").Value = "09"
$ws.Range("C2").Value = "Q2"

$ws.Range("A3").Value = "2023"
$ws.Range("B3").Value = "12"
$ws.Range("C3").Value = "Q3"

$ws.Range("A4").Value = "2024"
$ws.Range("B4").Value = "03"
$ws.Range("C4").Value = "Q4"

$ws.Range("A5").Value = "2024"
$ws.Range("B5").Value = "06"
$ws.Range("C5").Value = "Q1"

$ws.Range("A6").Value = "2024"
$ws.Range("B6").Value = "09"
$ws.Range("C6").Value = "Q2"
